$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original font formatting used in columns C and D so it can be
# re-applied after the hyperlinks are rebuilt (inserting a hyperlink resets a
# cell to the built-in "Hyperlink" style otherwise).
$cFontName = $ws.Range("C6").Font.Name
$cFontSize = $ws.Range("C6").Font.Size
$cFontColor = $ws.Range("C6").Font.Color
$cFontUnderline = $ws.Range("C6").Font.Underline

$dFontName = $ws.Range("D6").Font.Name
$dFontSize = $ws.Range("D6").Font.Size
$dFontColor = $ws.Range("D6").Font.Color
$dFontUnderline = $ws.Range("D6").Font.Underline

# Row 4 now holds a different reviewer's data: same app, but the keyword,
# date, reviewer emails and review comment all change.
$ws.Range("B4").Value = "helix jump"
$ws.Range("C4").Value = "nevilgreen@gmail.com"
$ws.Range("D4").Value = "vikicrestina@gmail.com"
$ws.Range("E4").Value = "27/5/2019 15:55"
$ws.Range("F4").Value = "good times with this game app. Helix jump is a wonderful game."

# Rebuild the mailto hyperlinks: rows 2-3 keep pointing at the same
# addresses, row 4 now points at the new reviewer's addresses. Rows 5-6 no
# longer have any hyperlinks once their data is cleared below.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:leviadlevi22@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "leviadlevi22@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:gazittalia1@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "gazittalia1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:sm6502345@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "sm6502345@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:cybworking@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "cybworking@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:nevilgreen@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "nevilgreen@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:vikicrestina@gmail.com", [System.Type]::Missing, [System.Type]::Missing, "vikicrestina@gmail.com")

foreach ($cellRef in @("C2", "C3", "C4")) {
    $ws.Range($cellRef).Font.Name = $cFontName
    $ws.Range($cellRef).Font.Size = $cFontSize
    $ws.Range($cellRef).Font.Color = $cFontColor
    $ws.Range($cellRef).Font.Underline = $cFontUnderline
}
foreach ($cellRef in @("D2", "D3", "D4")) {
    $ws.Range($cellRef).Font.Name = $dFontName
    $ws.Range($cellRef).Font.Size = $dFontSize
    $ws.Range($cellRef).Font.Color = $dFontColor
    $ws.Range($cellRef).Font.Underline = $dFontUnderline
}

# The two extra reviews that used to live in rows 5 and 6 are gone now.
$ws.Range("A5:F6").ClearContents()

# Clearing the rows nudges their auto height up slightly.
$ws.Rows.Item(5).RowHeight = 13.8
$ws.Rows.Item(6).RowHeight = 13.8

$ws.Range("F4").Select()
